# Update countries & provincias Spain
# Refresh the "paises" COVID-19 stats table with newer figures, swap the
# Etiopia / Bielorrusia rows (re-sorted by total cases) and bump the
# "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados" timestamp (A1) -----------------------------------
$ws.Range("A1").Value = "Datos actualizados a 9 de Octubre de 2020 a las 23:24"

# --- Row 4: Estados Unidos --------------------------------------------------
$ws.Range("B4").Value = 7879565
$ws.Range("C4").Value = 45802
$ws.Range("D4").Value = 5050985
$ws.Range("E4").Value = 2610224
$ws.Range("G4").Value = 618
$ws.Range("H4").Value = 218356

# --- Row 6: Brasil -----------------------------------------------------------
$ws.Range("B6").Value = 5055888
$ws.Range("C6").Value = 26349
$ws.Range("E6").Value = 491685
$ws.Range("G6").Value = 605
$ws.Range("H6").Value = 149639

# --- Row 25: Alemania ---------------------------------------------------------
$ws.Range("B25").Value = 320478
$ws.Range("C25").Value = 4964
$ws.Range("E25").Value = 41291
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 9687

# --- Row 29: Canada ------------------------------------------------------------
$ws.Range("B29").Value = 177719
$ws.Range("C29").Value = 2160
$ws.Range("D29").Value = 149270
$ws.Range("E29").Value = 18863

# --- Rows 53/54: Etiopia and Bielorrusia swap places (re-sorted) + refresh ----
$ws.Range("A53").Value = "Etiopia"
$ws.Range("B53").Value = 82662
$ws.Range("C53").Value = 865
$ws.Range("D53").Value = 37102
$ws.Range("E53").Value = 44289
$ws.Range("G53").Value = 9
$ws.Range("H53").Value = 1271

$ws.Range("A54").Value = "Bielorrusia"
$ws.Range("B54").Value = 82471
$ws.Range("C54").Value = 489
$ws.Range("D54").Value = 76543
$ws.Range("E54").Value = 5043
$ws.Range("G54").Value = 5
$ws.Range("H54").Value = 885

# --- Row 57: Barein --------------------------------------------------------
$ws.Range("B57").Value = 74860
$ws.Range("C57").Value = 438
$ws.Range("D57").Value = 70406
$ws.Range("E57").Value = 4183
$ws.Range("G57").Value = 7
$ws.Range("H57").Value = 271

# --- Row 90: Costa de Marfil -------------------------------------------------
$ws.Range("B90").Value = 20036
$ws.Range("C90").Value = 54
$ws.Range("D90").Value = 19696
$ws.Range("E90").Value = 220

# --- Row 118: Cabo Verde ------------------------------------------------------
$ws.Range("B118").Value = 6809
$ws.Range("C118").Value = 92
$ws.Range("D118").Value = 5932
$ws.Range("E118").Value = 804
$ws.Range("G118").Value = 2
$ws.Range("H118").Value = 73

# --- Row 132: Ruanda -----------------------------------------------------------
$ws.Range("B132").Value = 4890
$ws.Range("C132").Value = 5
$ws.Range("D132").Value = 3555
$ws.Range("E132").Value = 1305
$ws.Range("G132").Value = 1
$ws.Range("H132").Value = 30

# --- Row 138: Aruba -------------------------------------------------------------
$ws.Range("B138").Value = 4167
$ws.Range("C138").Value = 17
$ws.Range("E138").Value = 417
$ws.Range("G138").Value = 1
$ws.Range("H138").Value = 32

# --- Row 146: Guyana -------------------------------------------------------------
$ws.Range("B146").Value = 3358
$ws.Range("C146").Value = 29
$ws.Range("D146").Value = 2180
$ws.Range("E146").Value = 1078
$ws.Range("G146").Value = 2
$ws.Range("H146").Value = 100

# --- Row 194: Seychelles -----------------------------------------------------
$ws.Range("D194").Value = 144
$ws.Range("E194").Value = 4
